$wb = $excel.ActiveWorkbook

# The existing "总计" sheet is currently the 3rd sheet; insert the new
# "2022-Q1" sheet immediately before it so the final tab order becomes:
# 2021-Q1, 2021-Q4, 2022-Q1, 总计
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item(3))
$q1.Name = "2022-Q1"

# NOTE: Worksheets.Item(3) is a positional lookup, not a stable handle —
# once Add() inserts the new sheet at position 3, the old "总计" sheet has
# shifted to position 4. Re-resolve it now, after the insert.
$totalSheet = $wb.Worksheets.Item(4)

# Text-valued columns (B..G) must stay text (preserves leading zeros in
# fund codes and matches the source data's text-typed numeric-looking
# strings); force that with a text number format before assigning values.
$q1.Range("B1:G6").NumberFormat = "@"

# Headers (row 1)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "160916"
$q1.Range("C2").Value = "大成优选混合(LOF)"
$q1.Range("D2").Value = "16.14"
$q1.Range("E2").Value = "89.35"
$q1.Range("F2").Value = "7.27"
$q1.Range("G2").Value = "1.1734"
$q1.Range("H2").Value = 3

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "010738"
$q1.Range("C3").Value = "大成优选升级一年持有期混合A"
$q1.Range("D3").Value = "3.79"
$q1.Range("E3").Value = "89.02"
$q1.Range("F3").Value = "7.76"
$q1.Range("G3").Value = "0.2941"
$q1.Range("H3").Value = 5

# Row 4
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "007518"
$q1.Range("C4").Value = "东方阿尔法优选混合A"
$q1.Range("D4").Value = "2.03"
$q1.Range("E4").Value = "72.64"
$q1.Range("F4").Value = "1.40"
$q1.Range("G4").Value = "0.0284"
$q1.Range("H4").Value = 9

# Row 5
$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "007519"
$q1.Range("C5").Value = "东方阿尔法优选混合C"
$q1.Range("D5").Value = "0.82"
$q1.Range("E5").Value = "72.64"
$q1.Range("F5").Value = "1.40"
$q1.Range("G5").Value = "0.0115"
$q1.Range("H5").Value = 9

# Row 6
$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "010739"
$q1.Range("C6").Value = "大成优选升级一年持有期混合C"
$q1.Range("D6").Value = "0.09"
$q1.Range("E6").Value = "89.02"
$q1.Range("F6").Value = "7.76"
$q1.Range("G6").Value = "0.0070"
$q1.Range("H6").Value = 5

# Match the look of the sibling sheets: bold/centered/bordered header row
# and bold/centered/bordered index column (style used throughout the other
# quarter sheets). Copying formats from "2021-Q4" (already styled this way)
# keeps the exact same style instead of inventing a new one.
$styleSource = $wb.Worksheets.Item(2)
$styleSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Range("A2:A4").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now update the "总计" (totals) sheet: insert a new leading data row for
# 2022-Q1 (5 funds held, 1.51 亿元 held value) ahead of the existing rows,
# pushing 2021-Q4 and 2021-Q1 down by one row each, and renumber the index
# column (A) accordingly.
#
# Row 4 is brand new (the sheet only had rows 1-3 before), so first clone
# the index-column style from row 3 onto it, then fill in the values.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.45

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.76

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 1.51

# Restore the originally active sheet/tab (adding a worksheet makes it the
# active one, which the source diff does not change).
$wb.Worksheets.Item(1).Activate()
